# Update EIA Table 2.3.C workbook from the "October 2016" monthly release to
# the "November 2016" monthly release.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Update title text on row 2 (October -> November)
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "by Sector, 2006-November 2016 (Thousand Tons)"

# ---------------------------------------------------------------------
# 2. Insert a new row for "November" 2016 monthly data right after the
#    existing "October" row (row 52), pushing everything below down by one.
# ---------------------------------------------------------------------
$ws.Rows(53).Insert()

# Copy the formatting (styles/number formats/borders) from the row above
# (October, row 52) down onto the freshly inserted blank row 53, so the new
# row matches the existing month rows exactly.
$ws.Range("A52:F52").Copy()
$ws.Range("A53:F53").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new "November" row with its data.
$ws.Range("A53").Value = "November"
$ws.Range("B53").Value = 381
$ws.Range("C53").Value = 240
$ws.Range("D53").Value = 56
$ws.Range("E53").Value = 0.23
$ws.Range("F53").Value = 85

# ---------------------------------------------------------------------
# 3. Refresh the "Year to Date" figures (now rows 55-57 after the insert).
# ---------------------------------------------------------------------
$ws.Range("A55").Value = 2014
$ws.Range("B55").Value = 5157
$ws.Range("C55").Value = 3108
$ws.Range("D55").Value = 619
$ws.Range("E55").Value = 16
$ws.Range("F55").Value = 1413

$ws.Range("A56").Value = 2015
$ws.Range("B56").Value = 4826
$ws.Range("C56").Value = 2896
$ws.Range("D56").Value = 742
$ws.Range("E56").Value = 17
$ws.Range("F56").Value = 1171

$ws.Range("A57").Value = 2016
$ws.Range("B57").Value = 4862
$ws.Range("C57").Value = 3172
$ws.Range("D57").Value = 639
$ws.Range("E57").Value = 9
$ws.Range("F57").Value = 1042

# ---------------------------------------------------------------------
# 4. Update the "Rolling 12 Months" header text and figures
#    (now rows 58-60 after the insert).
# ---------------------------------------------------------------------
$ws.Range("A58").Value = "Rolling 12 Months Ending in November"

$ws.Range("A59").Value = 2015
$ws.Range("B59").Value = 5364
$ws.Range("C59").Value = 3232
$ws.Range("D59").Value = 811
$ws.Range("E59").Value = 19
$ws.Range("F59").Value = 1303

$ws.Range("A60").Value = 2016
$ws.Range("B60").Value = 5224
$ws.Range("C60").Value = 3404
$ws.Range("D60").Value = 676
$ws.Range("E60").Value = 10
$ws.Range("F60").Value = 1135
